$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.740.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.185.70'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.01%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.95%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.99%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.612'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.25%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.25%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.196.65'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.66%  '

$ws.Range("E10").Value = '  -3.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.84'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.41%  '

$ws.Range("E12").Value = '  -2.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.738.87'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.27%  '

$ws.Range("E14").Value = '  -2.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.742.12'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.19%  '

$ws.Range("E16").Value = '  -2.63%  '

$ws.Range("E17").Value = '  -3.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.197.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '420.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.01%  '

$ws.Range("E21").Value = '  -3.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.79%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("E24").Value = '  -1.84%  '

$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("E26").Value = '  +1.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.501'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.57%  '

$ws.Range("E28").Value = '  -6.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.89'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("E31").Value = '  -4.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.87'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.40%  '

$ws.Range("E33").Value = '  -0.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.63%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '156.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.38'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.95%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.707.03'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.10%  '

$ws.Range("E41").Value = '  -1.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '24.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.41%  '

$ws.Range("E43").Value = '  -5.02%  '

$ws.Range("E44").Value = '  -1.61%  '

$ws.Range("E45").Value = '  -5.44%  '

$ws.Range("E46").Value = '  -4.82%  '

$ws.Range("E47").Value = '  -3.22%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '294.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.91%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -11.56%  '

$ws.Range("E51").Value = '  -6.21%  '
